# QA renew By Richard 20160114
# Fill in answers that arrived for the still-open questions on the "QA"
# sheet, update the Logo answer with extra detail, and drop the two
# trailing blank rows that are no longer part of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QA")

# Row 4 - Logo answer gets more detail (AI file, resizable).
$ws.Range("E4").Value = "已提供，但大小像素是否合适正在确认中。提供的是AI文件，大小可以任意调整"

# Row 10 - Tab content question now answered.
$ws.Range("E10").Value = "就按照商品详情、消费提示、商家规则三个tab就好"
$ws.Range("F10").Value = "1/14/2016"
$ws.Range("G10").Value = "Lan"

# Row 11 - invoice PDF format question now answered.
$ws.Range("E11").Value = "这个我们会提供的"
$ws.Range("F11").Value = "1/14/2016"
$ws.Range("G11").Value = "Lan"

# Row 12 - invoice email question now answered.
$ws.Range("E12").Value = "需要填写，手机注册的用户注册后最好要强制完善个人信息，邮箱要填上去。或者按你说的在发invoice的时候填写。"
$ws.Range("F12").Value = "1/14/2016"
$ws.Range("G12").Value = "Lan"

# Row 13 - shipping / delivery address question now answered.
$ws.Range("E13").Value = "需要，并且应该有一个能自动计算运费的功能，这个澳大利亚邮局有api可以使用。货到付款的方式，我们自己有一个列表，里面列出了不同的社区不同的收费价格表，货到付款填写地址时 suburb是一个下拉菜单，只能在其中选择，其他的地方我们是不送的，对应的价格有一个excel表，到时一并发你"
$ws.Range("F13").Value = "1/14/2016"
$ws.Range("G13").Value = "Lan"

# Drop the two trailing empty rows (14, 15) - table now ends at row 13.
$ws.Rows.Item(15).Delete()
$ws.Rows.Item(14).Delete()
